# Update odds values on Sheet1 (rows 3, 6, 7) to the refreshed FlashScore data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("H3").Value = 3.4
$ws.Range("O3").Value = 1.25
$ws.Range("P3").Value = 3.75
$ws.Range("Q3").Value = 1.88
$ws.Range("R3").Value = 1.98
$ws.Range("AC3").Value = 11
$ws.Range("AE3").Value = 13
$ws.Range("AF3").Value = 41
$ws.Range("AI3").Value = 11
$ws.Range("AK3").Value = 23
$ws.Range("AU3").Value = 7.5
$ws.Range("BA3").Value = 67

# Row 6
$ws.Range("G6").Value = 1.21
$ws.Range("H6").Value = 5.6
$ws.Range("I6").Value = 10.5
$ws.Range("J6").Value = 1.6
$ws.Range("K6").Value = 2.72
$ws.Range("L6").Value = 8.25
$ws.Range("N6").Value = 10
$ws.Range("O6").Value = 1.13
$ws.Range("P6").Value = 5.2
$ws.Range("Q6").Value = 1.4
$ws.Range("R6").Value = 2.72
$ws.Range("S6").Value = 1.25
$ws.Range("T6").Value = 3.55
$ws.Range("U6").Value = 1.87
$ws.Range("V6").Value = 1.83
$ws.Range("W6").Value = 9.75
$ws.Range("Y6").Value = 9.25
$ws.Range("Z6").Value = 7.8
$ws.Range("AB6").Value = 25
$ws.Range("AD6").Value = 12.5
$ws.Range("AE6").Value = 23
$ws.Range("AF6").Value = 90
$ws.Range("AG6").Value = 35
$ws.Range("AH6").Value = 90
$ws.Range("AI6").Value = 35
$ws.Range("AJ6").Value = 350
$ws.Range("AK6").Value = 120
$ws.Range("AL6").Value = 90
$ws.Range("AM6").Value = 600
$ws.Range("AO6").Value = 5.2
$ws.Range("AQ6").Value = 12
$ws.Range("AT6").Value = 3.55
$ws.Range("AU6").Value = 9
$ws.Range("AV6").Value = 75
$ws.Range("AW6").Value = 11
$ws.Range("AX6").Value = 60
$ws.Range("AY6").Value = 50
$ws.Range("BA6").Value = 400

# Row 7
$ws.Range("G7").Value = 5.6
$ws.Range("I7").Value = 1.52
$ws.Range("K7").Value = 2.2
$ws.Range("N7").Value = 7.4
$ws.Range("P7").Value = 3.25
$ws.Range("Q7").Value = 1.9
$ws.Range("U7").Value = 1.98
$ws.Range("W7").Value = 14.5
$ws.Range("X7").Value = 35
$ws.Range("Y7").Value = 18.5
$ws.Range("Z7").Value = 120
$ws.Range("AC7").Value = 7.4
$ws.Range("AF7").Value = 90
$ws.Range("AG7").Value = 6.2
$ws.Range("AH7").Value = 6.7
$ws.Range("AJ7").Value = 10.5
$ws.Range("AM7").Value = 800
$ws.Range("AN7").Value = 7.1
$ws.Range("AP7").Value = 37
$ws.Range("AR7").Value = 250
$ws.Range("AV7").Value = 80
$ws.Range("AY7").Value = 18.5
